$wb = $excel.ActiveWorkbook

# Sheet that is currently "hotel_info" (index 1) will become "review_info".
# Sheet that is currently "review_info" (index 2) will become "hotel_info"
# (with a new "State" column). The underlying sheetId/r:id for each physical
# sheet stay put; only the name + contents are swapped, which is what the
# target workbook.xml / sheetN.xml actually reflect.
$shHotel  = $wb.Worksheets.Item(1)
$shReview = $wb.Worksheets.Item(2)

function Set-TextCell($cell, $text) {
    # Force the cell to be stored as a text/string value even when the
    # text looks like a number (e.g. "1166"), matching the source data
    # where these columns are shared strings, not numerics.
    $cell.NumberFormat = "@"
    $cell.Value2 = $text
}

# ---- Rebuild the sheet that will become "review_info" ----
$shHotel.Cells.Clear()

$reviewHeaders = @(
    "STR",
    "reviewer_ID",
    "reviewer_name",
    "Review_ID",
    "Date_of_scraping",
    "ReviewURL",
    "Tripadvisor_gcode",
    "Tripadvisor_dcode",
    "Tripadvisor_rcode",
    "review_date",
    "review_title",
    "review_content",
    "review_rating",
    "trip_month",
    "trip_purpose",
    "value",
    "rooms",
    "Location",
    "Cleanliness",
    "Sleep Quality",
    "Service",
    "Picture(yes=1)",
    "respondent",
    "response_date",
    "response_text"
)

for ($i = 0; $i -lt $reviewHeaders.Length; $i++) {
    $shHotel.Cells.Item(1, $i + 1).Value2 = $reviewHeaders[$i]
}

$shHotel.Name = "review_info_tmp"

# ---- Rebuild the sheet that will become "hotel_info" (adds "State") ----
$shReview.Cells.Clear()

$hotelHeaders = @(
    "STR",
    "Hotel_Name",
    "State",
    "City",
    "Zip",
    "TA_ReviewURL",
    "Tripadvisor_Hotel_Name",
    "English_Reviews_num",
    "Local_Rank",
    "Total_Reviews_num"
)

for ($i = 0; $i -lt $hotelHeaders.Length; $i++) {
    $shReview.Cells.Item(1, $i + 1).Value2 = $hotelHeaders[$i]
}

$shReview.Cells.Item(2, 1).Value2 = 4351
$shReview.Cells.Item(2, 2).Value2 = "Holiday Inn New Orleans Downtown Superdome"
$shReview.Cells.Item(2, 3).Value2 = "Louisiana"
$shReview.Cells.Item(2, 4).Value2 = "New Orleans"
$shReview.Cells.Item(2, 5).Value2 = 70112
$shReview.Cells.Item(2, 6).Value2 = "https://www.tripadvisor.com/Hotel_Review-g60864-d223121-Reviews-Holiday_Inn_Downtown_Superdome-New_Orleans_Louisiana.html"
$shReview.Cells.Item(2, 7).Value2 = "Holiday Inn Downtown Superdome"
Set-TextCell $shReview.Cells.Item(2, 8) "1166"
Set-TextCell $shReview.Cells.Item(2, 9) "95"
Set-TextCell $shReview.Cells.Item(2, 10) "1309"

$shReview.Name = "hotel_info"
$shHotel.Name = "review_info"

# ---- Put "review_info" tab before "hotel_info" tab ----
$reviewSheet = $wb.Worksheets.Item("review_info")
$hotelSheet  = $wb.Worksheets.Item("hotel_info")
$reviewSheet.Move($wb.Worksheets.Item(1))

$wb.Worksheets.Item(1).Activate()
